$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 110.85714  # H4: 264.33334 -> 110.85714
$ws.Cells.Item(4, 9).Value = 110.85714  # I4: 264.33334 -> 110.85714
$ws.Cells.Item(4, 11).Value = 110.85714  # K4: 264.33334 -> 110.85714
$ws.Cells.Item(4, 13).Value = 3.142859999999999  # M4: -150.33334 -> 3.142859999999999

# Row 43
$ws.Cells.Item(43, 8).Value = 46167.59  # H43: 375999.5 -> 46167.59
$ws.Cells.Item(43, 9).Value = 1958.25  # I43: 1999 -> 1958.25
$ws.Cells.Item(43, 10).Value = 152270  # J43: 750000 -> 152270
$ws.Cells.Item(43, 11).Value = 1958.25  # K43: 1999 -> 1958.25
$ws.Cells.Item(43, 12).Value = 152270  # L43: 750000 -> 152270
$ws.Cells.Item(43, 13).Value = -1889.25  # M43: -1930 -> -1889.25
$ws.Cells.Item(43, 14).Value = -152408  # N43: -750138 -> -152408

# Row 106
$ws.Cells.Item(106, 8).Value = 4976  # H106: 3634 -> 4976
$ws.Cells.Item(106, 9).Value = 4976  # I106: 3634 -> 4976
$ws.Cells.Item(106, 11).Value = 4976  # K106: 3634 -> 4976
$ws.Cells.Item(106, 13).Value = -4345  # M106: -3003 -> -4345

# Row 113
$ws.Cells.Item(113, 8).Value = 4555  # H113: 4641.4287 -> 4555
$ws.Cells.Item(113, 9).Value = 4305  # I113: 4482.5 -> 4305
$ws.Cells.Item(113, 11).Value = 4305  # K113: 4482.5 -> 4305
$ws.Cells.Item(113, 13).Value = -1051  # M113: -1228.5 -> -1051

# Row 127
$ws.Cells.Item(127, 8).Value = 2354.8572  # H127: 2355.7144 -> 2354.8572
$ws.Cells.Item(127, 9).Value = 1580.6666  # I127: 1581.6666 -> 1580.6666
$ws.Cells.Item(127, 11).Value = 4741.9998  # K127: 4744.9998 -> 4741.9998
$ws.Cells.Item(127, 13).Value = 218.0002000000004  # M127: 215.0002000000004 -> 218.0002000000004

# Row 137
$ws.Cells.Item(137, 8).Value = 1307.5883  # H137: 1408.8 -> 1307.5883
$ws.Cells.Item(137, 9).Value = 1146.625  # I137: 1232.0714 -> 1146.625
$ws.Cells.Item(137, 11).Value = 3439.875  # K137: 3696.2142 -> 3439.875
$ws.Cells.Item(137, 13).Value = -889.875  # M137: -1146.2142 -> -889.875

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 8091.905  # H32: 8862.157999999999 -> 8091.905
$ws.Cells.Item(32, 9).Value = 7320.5264  # I32: 8090.647 -> 7320.5264
$ws.Cells.Item(32, 11).Value = 7320.5264  # K32: 8090.647 -> 7320.5264
$ws.Cells.Item(32, 13).Value = -7033.5264  # M32: -7803.647 -> -7033.5264

# Row 37
$ws.Cells.Item(37, 8).Value = 19000  # H37: 17500 -> 19000

# Row 102
$ws.Cells.Item(102, 8).Value = 2738.9285  # H102: 3112.3333 -> 2738.9285
$ws.Cells.Item(102, 9).Value = 2572.2727  # I102: 3033.111 -> 2572.2727
$ws.Cells.Item(102, 11).Value = 2572.2727  # K102: 3033.111 -> 2572.2727
$ws.Cells.Item(102, 13).Value = -950.2727  # M102: -1411.111 -> -950.2727

# Row 122
$ws.Cells.Item(122, 8).Value = 2998  # H122: 2997 -> 2998
$ws.Cells.Item(122, 10).Value = 3000  # J122: 0 -> 3000
$ws.Cells.Item(122, 12).Value = 9000  # L122: 0 -> 9000
$ws.Cells.Item(122, 14).Value = -13900  # N122: None -> -13900

# Row 132
$ws.Cells.Item(132, 8).Value = 1951.5  # H132: 2624.25 -> 1951.5
$ws.Cells.Item(132, 9).Value = 1941.8  # I132: 2832.3333 -> 1941.8
$ws.Cells.Item(132, 11).Value = 5825.4  # K132: 8496.999899999999 -> 5825.4
$ws.Cells.Item(132, 13).Value = -3295.4  # M132: -5966.999899999999 -> -3295.4

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Cells.Item(134, 8).Value = 8103  # H134: 8972.9 -> 8103
$ws.Cells.Item(134, 9).Value = 8293.444  # I134: 8866.125 -> 8293.444
$ws.Cells.Item(134, 10).Value = 7531.6665  # J134: 9400 -> 7531.6665
$ws.Cells.Item(134, 11).Value = 24880.332  # K134: 26598.375 -> 24880.332
$ws.Cells.Item(134, 12).Value = 22594.9995  # L134: 28200 -> 22594.9995
$ws.Cells.Item(134, 13).Value = -22345.332  # M134: -24063.375 -> -22345.332
$ws.Cells.Item(134, 14).Value = -27664.9995  # N134: -33270 -> -27664.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 1349.8  # H5: 822.53845 -> 1349.8
$ws.Cells.Item(5, 9).Value = 583  # I5: 517.5454999999999 -> 583
$ws.Cells.Item(5, 11).Value = 1749  # K5: 1552.6365 -> 1749
$ws.Cells.Item(5, 13).Value = -1637  # M5: -1440.6365 -> -1637

# Row 26
$ws.Cells.Item(26, 8).Value = 125692.5  # H26: 262775 -> 125692.5
$ws.Cells.Item(26, 9).Value = 333510  # I26: 500250 -> 333510
$ws.Cells.Item(26, 10).Value = 1002  # J26: 25300 -> 1002
$ws.Cells.Item(26, 11).Value = 1000530  # K26: 1500750 -> 1000530
$ws.Cells.Item(26, 12).Value = 3006  # L26: 75900 -> 3006
$ws.Cells.Item(26, 13).Value = -1000242  # M26: -1500462 -> -1000242
$ws.Cells.Item(26, 14).Value = -3582  # N26: -76476 -> -3582

# Row 74
$ws.Cells.Item(74, 8).Value = 500000  # H74: 0 -> 500000
$ws.Cells.Item(74, 10).Value = 500000  # J74: 0 -> 500000
$ws.Cells.Item(74, 12).Value = 1500000  # L74: 0 -> 1500000
$ws.Cells.Item(74, 14).Value = -1502122  # N74: None -> -1502122

# Row 77
$ws.Cells.Item(77, 8).Value = 500000  # H77: 0 -> 500000
$ws.Cells.Item(77, 10).Value = 500000  # J77: 0 -> 500000
$ws.Cells.Item(77, 12).Value = 4500000  # L77: 0 -> 4500000
$ws.Cells.Item(77, 14).Value = -4510608  # N77: None -> -4510608

# Row 135
$ws.Cells.Item(135, 8).Value = 1349.8  # H135: 822.53845 -> 1349.8
$ws.Cells.Item(135, 9).Value = 583  # I135: 517.5454999999999 -> 583
$ws.Cells.Item(135, 11).Value = 5247  # K135: 4657.9095 -> 5247
$ws.Cells.Item(135, 13).Value = -2712  # M135: -2122.9095 -> -2712

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Cells.Item(11, 8).Value = 10000000  # H11: 15002500 -> 10000000
$ws.Cells.Item(11, 9).Value = 0  # I11: 70000000 -> 0
$ws.Cells.Item(11, 10).Value = 10000000  # J11: 7145715 -> 10000000
$ws.Cells.Item(11, 11).Value = 0  # K11: 70000000 -> 0
$ws.Cells.Item(11, 12).Value = 10000000  # L11: 7145715 -> 10000000
$ws.Cells.Item(11, 13).ClearContents()  # M11: -69999861 -> (removed)
$ws.Cells.Item(11, 14).Value = -10000278  # N11: -7145993 -> -10000278

# Row 41
$ws.Cells.Item(41, 8).Value = 2391.6667  # H41: 4440 -> 2391.6667
$ws.Cells.Item(41, 9).Value = 4350  # I41: 3300 -> 4350
$ws.Cells.Item(41, 10).Value = 2000  # J41: 9000 -> 2000
$ws.Cells.Item(41, 11).Value = 4350  # K41: 3300 -> 4350
$ws.Cells.Item(41, 12).Value = 2000  # L41: 9000 -> 2000
$ws.Cells.Item(41, 13).Value = -3995  # M41: -2945 -> -3995
$ws.Cells.Item(41, 14).Value = -2710  # N41: -9710 -> -2710

# Row 80
$ws.Cells.Item(80, 8).Value = 1713.0435  # H80: 1707.96 -> 1713.0435
$ws.Cells.Item(80, 9).Value = 1777.6666  # I80: 1726.4 -> 1777.6666
$ws.Cells.Item(80, 11).Value = 1777.6666  # K80: 1726.4 -> 1777.6666
$ws.Cells.Item(80, 13).Value = -779.6666  # M80: -728.4000000000001 -> -779.6666

# Row 83
$ws.Cells.Item(83, 8).Value = 1713.0435  # H83: 1707.96 -> 1713.0435
$ws.Cells.Item(83, 9).Value = 1777.6666  # I83: 1726.4 -> 1777.6666
$ws.Cells.Item(83, 11).Value = 8888.333000000001  # K83: 8632 -> 8888.333000000001
$ws.Cells.Item(83, 13).Value = -3896.333000000001  # M83: -3640 -> -3896.333000000001

# Row 101
$ws.Cells.Item(101, 8).Value = 31871.334  # H101: 13057 -> 31871.334
$ws.Cells.Item(101, 10).Value = 31871.334  # J101: 13057 -> 31871.334
$ws.Cells.Item(101, 12).Value = 31871.334  # L101: 13057 -> 31871.334
$ws.Cells.Item(101, 14).Value = -38361.334  # N101: -19547 -> -38361.334

# Row 122
$ws.Cells.Item(122, 8).Value = 10420029  # H122: 11367213 -> 10420029
$ws.Cells.Item(122, 9).Value = 12503235  # I122: 13892372 -> 12503235
$ws.Cells.Item(122, 11).Value = 37509705  # K122: 41677116 -> 37509705
$ws.Cells.Item(122, 13).Value = -37507255  # M122: -41674666 -> -37507255

# Row 132
$ws.Cells.Item(132, 8).Value = 1490.6666  # H132: 1618 -> 1490.6666
$ws.Cells.Item(132, 10).Value = 0  # J132: 2000 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 6000 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -11060 -> (removed)

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 3286.2856  # H7: 2000 -> 3286.2856
$ws.Cells.Item(7, 9).Value = 2251  # I7: 2000 -> 2251
$ws.Cells.Item(7, 10).Value = 4666.6665  # J7: 0 -> 4666.6665
$ws.Cells.Item(7, 11).Value = 2251  # K7: 2000 -> 2251
$ws.Cells.Item(7, 12).Value = 4666.6665  # L7: 0 -> 4666.6665
$ws.Cells.Item(7, 13).Value = -2139  # M7: -1888 -> -2139
$ws.Cells.Item(7, 14).Value = -4890.6665  # N7: None -> -4890.6665

# Row 12
$ws.Cells.Item(12, 8).Value = 2951  # H12: 2749.75 -> 2951
$ws.Cells.Item(12, 9).Value = 903  # I12: 1000 -> 903
$ws.Cells.Item(12, 10).Value = 4999  # J12: 3333 -> 4999
$ws.Cells.Item(12, 11).Value = 903  # K12: 1000 -> 903
$ws.Cells.Item(12, 12).Value = 4999  # L12: 3333 -> 4999
$ws.Cells.Item(12, 13).Value = -733  # M12: -830 -> -733
$ws.Cells.Item(12, 14).Value = -5339  # N12: -3673 -> -5339

# Row 29
$ws.Cells.Item(29, 8).Value = 10000  # H29: 28800 -> 10000
$ws.Cells.Item(29, 9).Value = 10000  # I29: 0 -> 10000
$ws.Cells.Item(29, 10).Value = 0  # J29: 28800 -> 0
$ws.Cells.Item(29, 11).Value = 10000  # K29: 0 -> 10000
$ws.Cells.Item(29, 12).Value = 0  # L29: 28800 -> 0
$ws.Cells.Item(29, 13).Value = -9705  # M29: None -> -9705
$ws.Cells.Item(29, 14).ClearContents()  # N29: -29390 -> (removed)

# Row 61
$ws.Cells.Item(61, 8).Value = 3836.8823  # H61: 3941.8667 -> 3836.8823
$ws.Cells.Item(61, 9).Value = 3681.4666  # I61: 3865.9285 -> 3681.4666
$ws.Cells.Item(61, 10).Value = 5002.5  # J61: 5005 -> 5002.5
$ws.Cells.Item(61, 11).Value = 3681.4666  # K61: 3865.9285 -> 3681.4666
$ws.Cells.Item(61, 12).Value = 5002.5  # L61: 5005 -> 5002.5
$ws.Cells.Item(61, 13).Value = -3479.4666  # M61: -3663.9285 -> -3479.4666
$ws.Cells.Item(61, 14).Value = -5406.5  # N61: -5409 -> -5406.5

# Row 113
$ws.Cells.Item(113, 8).Value = 3836.8823  # H113: 3941.8667 -> 3836.8823
$ws.Cells.Item(113, 9).Value = 3681.4666  # I113: 3865.9285 -> 3681.4666
$ws.Cells.Item(113, 10).Value = 5002.5  # J113: 5005 -> 5002.5
$ws.Cells.Item(113, 11).Value = 3681.4666  # K113: 3865.9285 -> 3681.4666
$ws.Cells.Item(113, 12).Value = 5002.5  # L113: 5005 -> 5002.5
$ws.Cells.Item(113, 13).Value = -1511.4666  # M113: -1695.9285 -> -1511.4666
$ws.Cells.Item(113, 14).Value = -9342.5  # N113: -9345 -> -9342.5

# Row 122
$ws.Cells.Item(122, 8).Value = 3650.75  # H122: 4126.25 -> 3650.75
$ws.Cells.Item(122, 9).Value = 3090  # I122: 3402.2 -> 3090
$ws.Cells.Item(122, 11).Value = 9270  # K122: 10206.6 -> 9270
$ws.Cells.Item(122, 13).Value = -6820  # M122: -7756.599999999999 -> -6820

# Row 126
$ws.Cells.Item(126, 8).Value = 3286.2856  # H126: 2000 -> 3286.2856
$ws.Cells.Item(126, 9).Value = 2251  # I126: 2000 -> 2251
$ws.Cells.Item(126, 10).Value = 4666.6665  # J126: 0 -> 4666.6665
$ws.Cells.Item(126, 11).Value = 6753  # K126: 6000 -> 6753
$ws.Cells.Item(126, 12).Value = 13999.9995  # L126: 0 -> 13999.9995
$ws.Cells.Item(126, 13).Value = -4283  # M126: -3530 -> -4283
$ws.Cells.Item(126, 14).Value = -18939.9995  # N126: None -> -18939.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 1653.5  # H122: 1864.3 -> 1653.5
$ws.Cells.Item(122, 9).Value = 1653.5  # I122: 1864.3 -> 1653.5
$ws.Cells.Item(122, 11).Value = 4960.5  # K122: 5592.9 -> 4960.5
$ws.Cells.Item(122, 13).Value = -2510.5  # M122: -3142.9 -> -2510.5

# Row 132
$ws.Cells.Item(132, 8).Value = 0  # H132: 945.9091 -> 0
$ws.Cells.Item(132, 9).Value = 0  # I132: 884.44446 -> 0
$ws.Cells.Item(132, 10).Value = 0  # J132: 1222.5 -> 0
$ws.Cells.Item(132, 11).Value = 0  # K132: 2653.33338 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 3667.5 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: -123.33338 -> (removed)
$ws.Cells.Item(132, 14).ClearContents()  # N132: -8727.5 -> (removed)

# Row 136
$ws.Cells.Item(136, 8).Value = 1406.1578  # H136: 1451 -> 1406.1578
$ws.Cells.Item(136, 9).Value = 1406.1578  # I136: 1451 -> 1406.1578
$ws.Cells.Item(136, 11).Value = 4218.4734  # K136: 4353 -> 4218.4734
$ws.Cells.Item(136, 13).Value = -1668.4734  # M136: -1803 -> -1668.4734
